$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("ICER") holds a sentinel value of 999999999 on a number of rows;
# cap those values down to 100000 so the ICER plots scale properly.
$rows = @(4, 18, 33, 47, 52, 57, 59, 63, 64, 66, 97, 105, 114, 133, 135, 136, 152, 155)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq 999999999) {
        $cell.Value = 100000
    }
}
